# Fix content placeholders in the slide master / slide layouts that were
# overlapping other template elements (footer / slide-number placeholders)
# near the bottom of the slide.
#
# NOTE: PowerPoint's COM object model exposes shape position/size in
# points, while the underlying OOXML stores EMU (1 pt = 12700 EMU). Some
# target EMU values are not an exact multiple of 12700, so a plain
# EMU/12700 conversion can round to the wrong integer EMU after the
# round-trip through the host. Adding a tiny epsilon (0.00005 pt, i.e.
# well under 1 EMU) nudges the conversion so it lands on the exact target
# EMU value every time, without perceptibly moving the shape.

$EMU_PER_POINT = 12700
$EPS = 0.00005

function ToPt([double]$emu) {
    return ($emu / $EMU_PER_POINT) + $EPS
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$layouts = $master.CustomLayouts

# -----------------------------------------------------------------
# Layout 2 - "Title and Content"
# -----------------------------------------------------------------
$lay2 = $layouts.Item(2)

# "Content Placeholder 2" (idx 1) previously had no explicit <p:spPr>
# xfrm (it inherited the oversized master placeholder). Give it its own,
# shorter, explicit size so its bottom no longer overlaps the footer row.
$contentPh = $lay2.Shapes.Item("Content Placeholder 2")
$contentPh.Left = ToPt 838200
$contentPh.Top = ToPt 1825625
$contentPh.Width = ToPt 10515600
$contentPh.Height = ToPt 3780045

# The redundant "Picture Placeholder 7" (idx 13) duplicated the same
# area as the content placeholder above - remove it entirely.
$picPh = $lay2.Shapes.Item("Picture Placeholder 7")
$picPh.Delete()

# -----------------------------------------------------------------
# Layout 3 - "Section Header"
# -----------------------------------------------------------------
$lay3 = $layouts.Item(3)

$title3 = $lay3.Shapes.Item("Title 1")
$title3.Top = ToPt 1242599

$text3 = $lay3.Shapes.Item("Text Placeholder 2")
$text3.Top = ToPt 4122324

# -----------------------------------------------------------------
# Layout 4 - "Two Content" (left placeholder only)
# -----------------------------------------------------------------
$lay4 = $layouts.Item(4)
$content4 = $lay4.Shapes.Item("Content Placeholder 2")
$content4.Height = ToPt 3809862

# -----------------------------------------------------------------
# Layout 5 - "Comparison" (left placeholder only)
# -----------------------------------------------------------------
$lay5 = $layouts.Item(5)
$content5 = $lay5.Shapes.Item("Content Placeholder 3")
$content5.Height = ToPt 3036891

# -----------------------------------------------------------------
# Layout 8 - "Content with Caption"
# -----------------------------------------------------------------
$lay8 = $layouts.Item(8)
$text8 = $lay8.Shapes.Item("Text Placeholder 3")
$text8.Height = ToPt 3518452

# -----------------------------------------------------------------
# Layout 9 - "Picture with Caption"
# -----------------------------------------------------------------
$lay9 = $layouts.Item(9)
$text9 = $lay9.Shapes.Item("Text Placeholder 3")
$text9.Height = ToPt 3518452

# -----------------------------------------------------------------
# Slide Master - body "Text Placeholder 2"
# -----------------------------------------------------------------
$masterText = $master.Shapes.Item("Text Placeholder 2")
$masterText.Height = ToPt 3832225

Write-Host "Template placeholder overlap fix applied."
